$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2-51: refresh crypto price/volume data (and ARBITRUM/HuobiToken rank swap at rows 33-34)
$ws.Range("D2").Value = "27.424.31"
$ws.Range("E2").Value = "  +2.19%  "

$ws.Range("D3").Value = "1.871.21"
$ws.Range("E3").Value = "  +3.10%  "

$ws.Range("E4").Value = "  -0.61%  "

$ws.Range("D5").Value = "315.42"
$ws.Range("E5").Value = "  +2.34%  "

$ws.Range("E6").Value = "  -0.54%  "

$ws.Range("D7").Value = "0.4671"
$ws.Range("E7").Value = "  +1.32%  "

$ws.Range("D8").Value = "0.3736"
$ws.Range("E8").Value = "  +2.84%  "

$ws.Range("E9").Value = "  +2.67%  "

$ws.Range("D10").Value = "0.8924"
$ws.Range("E10").Value = "  +4.18%  "

$ws.Range("D11").Value = "0.07956"
$ws.Range("E11").Value = "  +5.91%  "

$ws.Range("D12").Value = "20.18"
$ws.Range("E12").Value = "  +2.45%  "

$ws.Range("D13").Value = "1.858.94"
$ws.Range("E13").Value = "  +1.25%  "

$ws.Range("D14").Value = "5.454"
$ws.Range("E14").Value = "  +2.60%  "

$ws.Range("D15").Value = "6.628"
$ws.Range("E15").Value = "  +2.00%  "

$ws.Range("D16").Value = "92.99"
$ws.Range("E16").Value = "  +1.41%  "

$ws.Range("D17").Value = "1.003"
$ws.Range("E17").Value = "  -0.59%  "

$ws.Range("D18").Value = "0.000008983"
$ws.Range("E18").Value = "  +4.92%  "

$ws.Range("D19").Value = "1.003"
$ws.Range("E19").Value = "  -0.42%  "

$ws.Range("D20").Value = "14.98"
$ws.Range("E20").Value = "  +4.07%  "

$ws.Range("D21").Value = "27.441.25"
$ws.Range("E21").Value = "  +2.69%  "

$ws.Range("D22").Value = "5.182"
$ws.Range("E22").Value = "  +0.86%  "

$ws.Range("D23").Value = "10.64"
$ws.Range("E23").Value = "  +1.51%  "

$ws.Range("D24").Value = "2.070.95"
$ws.Range("E24").Value = "  +3.69%  "

$ws.Range("D25").Value = "152.76"
$ws.Range("E25").Value = "  +1.19%  "

$ws.Range("D26").Value = "1.879"
$ws.Range("E26").Value = "  +1.60%  "

$ws.Range("D27").Value = "18.62"
$ws.Range("E27").Value = "  +2.72%  "

$ws.Range("D28").Value = "2.107"
$ws.Range("E28").Value = "  +2.32%  "

$ws.Range("D29").Value = "5.191"
$ws.Range("E29").Value = "  +1.96%  "

$ws.Range("D30").Value = "117.66"
$ws.Range("E30").Value = "  +2.32%  "

$ws.Range("D31").Value = "0.08915"
$ws.Range("E31").Value = "  +0.71%  "

$ws.Range("D32").Value = "0.7578"
$ws.Range("E32").Value = "  +5.92%  "

$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").Value = "1.166"
$ws.Range("E33").Value = "  +3.37%  "

$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "2.968"
$ws.Range("E34").Value = "  +0.61%  "

$ws.Range("D35").Value = "4.518"
$ws.Range("E35").Value = "  +2.64%  "

$ws.Range("D36").Value = "2.669"
$ws.Range("E36").Value = "  +9.37%  "

$ws.Range("E37").Value = "  +1.18%  "

$ws.Range("E38").Value = "  +1.40%  "

$ws.Range("D39").Value = "0.01962"
$ws.Range("E39").Value = "  +2.96%  "

$ws.Range("D40").Value = "2.997"
$ws.Range("E40").Value = "  +2.67%  "

$ws.Range("D41").Value = "7.207"
$ws.Range("E41").Value = "  +0.87%  "

$ws.Range("D42").Value = "0.5249"
$ws.Range("E42").Value = "  +2.55%  "

$ws.Range("D43").Value = "0.1650"
$ws.Range("E43").Value = "  +1.98%  "

$ws.Range("D44").Value = "8.361"
$ws.Range("E44").Value = "  +2.33%  "

$ws.Range("D45").Value = "0.4935"
$ws.Range("E45").Value = "  +3.15%  "

$ws.Range("D46").Value = "10.32"
$ws.Range("E46").Value = "  +2.19%  "

$ws.Range("D48").Value = "103.81"
$ws.Range("E48").Value = "  +1.10%  "

$ws.Range("D49").Value = "1.658"
$ws.Range("E49").Value = "  +2.80%  "

$ws.Range("D50").Value = "0.06275"
$ws.Range("E50").Value = "  +1.49%  "

$ws.Range("D51").Value = "66.17"
$ws.Range("E51").Value = "  +3.50%  "
